# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
#
# The source feed re-sorted a handful of same-date fixtures, which shuffles
# which spreadsheet row each match's odds/result data lands on. Column A
# (the sequential "id" index) stays put; every other column (B..AD) for the
# affected rows is swapped/rotated to its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

# rows 100 <-> 101
Swap-Rows $ws 100 101

# rows 102 <-> 103
Swap-Rows $ws 102 103

# rows 162 <-> 163
Swap-Rows $ws 162 163

# rows 204 -> 205 -> 206 -> 204 (3-way cyclic rotation)
$range204 = $ws.Range("B204:AD204")
$range205 = $ws.Range("B205:AD205")
$range206 = $ws.Range("B206:AD206")

$vals204 = $range204.Value()
$vals205 = $range205.Value()
$vals206 = $range206.Value()

# after: 204 takes what was in 205, 205 takes what was in 206, 206 takes what was in 204
$range204.Value = $vals205
$range205.Value = $vals206
$range206.Value = $vals204
